# "Rimligare antal aktiva bilar"
# Adjust allotted times for a handful of individuals in the Förmiddag
# ('9-11') window, drop the "woman" constraint tag from Individ 40, and
# remove the trailing three rows (the extra Eftermiddag '13-15' entries
# for Individ 6, Individ 15 and Individ 18) that made the car count
# unreasonable.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Individ 2 ('Förmiddag', '9-11') - row 48
$ws.Range("B48").Value = 10

# Individ 6 ('Förmiddag', '9-11') - row 49
$ws.Range("B49").Value = 45
$ws.Range("D49").Value = " "

# Individ 6 ('Förmiddag', '9-11') - row 50 (duplicate entry)
$ws.Range("B50").Value = 45
$ws.Range("D50").Value = " "

# Individ 40 ('Förmiddag', '9-11') - row 52
$ws.Range("B52").Value = 10
$ws.Range("D52").Value = "license"

# Individ 53 ('Förmiddag', '9-11') - row 53
$ws.Range("B53").Value = 50

# Individ 61 ('Förmiddag', '9-11') - row 56
$ws.Range("B56").Value = 55

# Drop the trailing three rows (84-86: extra Individ 6, Individ 15,
# Individ 18 Eftermiddag '13-15' entries) - shrinks used range to A1:G83
$ws.Range("A84:G86").EntireRow.Delete()
